$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 132
$ws.Range("H132").Value = 23813738
$ws.Range("I132").Value = 29415900
$ws.Range("J132").Value = 4545.375
$ws.Range("K132").Value = 88247700
$ws.Range("L132").Value = 13636.125
$ws.Range("M132").Value = -88245170
$ws.Range("N132").Value = -18696.125

# Row 137
$ws.Range("H137").Value = 1986358
$ws.Range("I137").Value = 2802323
$ws.Range("J137").Value = 4728.5713
$ws.Range("K137").Value = 8406969
$ws.Range("L137").Value = 14185.7139
$ws.Range("M137").Value = -8404419
$ws.Range("N137").Value = -19285.7139

# Row 141
$ws.Range("H141").Value = 168512.67
$ws.Range("I141").Value = 287442.44
$ws.Range("J141").Value = 2011
$ws.Range("K141").Value = 862327.3200000001
$ws.Range("L141").Value = 6033
$ws.Range("M141").Value = -857147.3200000001
$ws.Range("N141").Value = -16393

$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 886.9286
$ws.Range("I2").Value = 810.8182
$ws.Range("K2").Value = 810.8182
$ws.Range("M2").Value = -697.8182

# Row 61
$ws.Range("H61").Value = 1954.8096
$ws.Range("I61").Value = 1182.091
$ws.Range("K61").Value = 1182.091
$ws.Range("M61").Value = -970.0909999999999

# Row 74
$ws.Range("H74").Value = 4563.25
$ws.Range("I74").Value = 5705.353
$ws.Range("K74").Value = 5705.353
$ws.Range("M74").Value = -4831.353

# Row 77
$ws.Range("H77").Value = 4563.25
$ws.Range("I77").Value = 5705.353
$ws.Range("K77").Value = 28526.765
$ws.Range("M77").Value = -24158.765

# Row 97
$ws.Range("H97").Value = 1826.4166
$ws.Range("I97").Value = 1090.6
$ws.Range("K97").Value = 1090.6
$ws.Range("M97").Value = -594.5999999999999

# Row 103
$ws.Range("H103").Value = 34191.65
$ws.Range("J103").Value = 34191.65
$ws.Range("L103").Value = 34191.65
$ws.Range("N103").Value = -36535.65

# Row 116
$ws.Range("H116").Value = 886.9286
$ws.Range("I116").Value = 810.8182
$ws.Range("K116").Value = 810.8182
$ws.Range("M116").Value = 1483.1818

# Row 136
$ws.Range("H136").Value = 1954.8096
$ws.Range("I136").Value = 1182.091
$ws.Range("K136").Value = 3546.273
$ws.Range("M136").Value = -996.2729999999997

$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 886.9286
$ws.Range("I3").Value = 810.8182
$ws.Range("K3").Value = 810.8182
$ws.Range("M3").Value = -696.8182

# Row 134
$ws.Range("H134").Value = 1937.5526
$ws.Range("I134").Value = 1271.7407
$ws.Range("J134").Value = 3571.818
$ws.Range("K134").Value = 3815.2221
$ws.Range("L134").Value = 10715.454
$ws.Range("M134").Value = -1280.2221
$ws.Range("N134").Value = -15785.454

$ws = $wb.Worksheets.Item("CRP")
# Row 58
$ws.Range("H58").Value = 2526.9177
$ws.Range("I58").Value = 1587.9642
$ws.Range("J58").Value = 5619.9414
$ws.Range("K58").Value = 1587.9642
$ws.Range("L58").Value = 5619.9414
$ws.Range("M58").Value = -1384.9642
$ws.Range("N58").Value = -6025.9414

# Row 63
$ws.Range("H63").Value = 49996
$ws.Range("J63").Value = 49996
$ws.Range("L63").Value = 49996
$ws.Range("N63").Value = -51368

# Row 66
$ws.Range("H66").Value = 49996
$ws.Range("J66").Value = 49996
$ws.Range("L66").Value = 149988
$ws.Range("N66").Value = -156852

# Row 100
$ws.Range("H100").Value = 68000
$ws.Range("J100").Value = 68000
$ws.Range("L100").Value = 68000
$ws.Range("N100").Value = -70164

# Row 132
$ws.Range("H132").Value = 2414
$ws.Range("I132").Value = 1374.2106
$ws.Range("J132").Value = 8999.333000000001
$ws.Range("K132").Value = 4122.6318
$ws.Range("L132").Value = 26997.999
$ws.Range("M132").Value = -1592.6318
$ws.Range("N132").Value = -32057.999

# Row 134
$ws.Range("H134").Value = 8923.933999999999
$ws.Range("I134").Value = 11146.1
$ws.Range("K134").Value = 33438.3
$ws.Range("M134").Value = -30903.3

# Row 136
$ws.Range("H136").Value = 2526.9177
$ws.Range("I136").Value = 1587.9642
$ws.Range("J136").Value = 5619.9414
$ws.Range("K136").Value = 4763.892599999999
$ws.Range("L136").Value = 16859.8242
$ws.Range("M136").Value = -2213.892599999999
$ws.Range("N136").Value = -21959.8242

# Row 137
$ws.Range("H137").Value = 40407.5
$ws.Range("J137").Value = 40407.5
$ws.Range("L137").Value = 40407.5
$ws.Range("N137").Value = -50607.5

$ws = $wb.Worksheets.Item("CUL")
# Row 113
$ws.Range("H113").Value = 3788475.5
$ws.Range("I113").Value = 610.2222
$ws.Range("J113").Value = 8333913.5
$ws.Range("K113").Value = 1830.6666
$ws.Range("L113").Value = 25001740.5
$ws.Range("M113").Value = 339.3334
$ws.Range("N113").Value = -25006080.5

$ws = $wb.Worksheets.Item("GSM")
# Row 70
$ws.Range("H70").Value = 6734.1904
$ws.Range("I70").Value = 5863.625
$ws.Range("K70").Value = 5863.625
$ws.Range("M70").Value = -5593.625

# Row 73
$ws.Range("H73").Value = 6734.1904
$ws.Range("I73").Value = 5863.625
$ws.Range("K73").Value = 5863.625
$ws.Range("M73").Value = -4927.625

# Row 126
$ws.Range("H126").Value = 3243.59
$ws.Range("I126").Value = 2880.329
$ws.Range("J126").Value = 4610.143
$ws.Range("K126").Value = 8640.987000000001
$ws.Range("L126").Value = 13830.429
$ws.Range("M126").Value = -6170.987000000001
$ws.Range("N126").Value = -18770.429

# Row 128
$ws.Range("H128").Value = 41816.668
$ws.Range("J128").Value = 41816.668
$ws.Range("L128").Value = 41816.668
$ws.Range("N128").Value = -51776.668

# Row 132
$ws.Range("H132").Value = 3137.2222
$ws.Range("I132").Value = 1837.9375
$ws.Range("J132").Value = 5027.091
$ws.Range("K132").Value = 5513.8125
$ws.Range("L132").Value = 15081.273
$ws.Range("M132").Value = -2983.8125
$ws.Range("N132").Value = -20141.273

$ws = $wb.Worksheets.Item("LTW")
# Row 132
$ws.Range("H132").Value = 5909.6816
$ws.Range("I132").Value = 1236.75
$ws.Range("J132").Value = 8579.929
$ws.Range("K132").Value = 3710.25
$ws.Range("L132").Value = 25739.787
$ws.Range("M132").Value = -1180.25
$ws.Range("N132").Value = -30799.787

# Row 136
$ws.Range("H136").Value = 4203.48
$ws.Range("I136").Value = 1186.6364
$ws.Range("J136").Value = 6573.857
$ws.Range("K136").Value = 3559.9092
$ws.Range("L136").Value = 19721.571
$ws.Range("M136").Value = -1009.9092
$ws.Range("N136").Value = -24821.571

$ws = $wb.Worksheets.Item("WVR")
# Row 107
$ws.Range("H107").Value = 792.2353000000001
$ws.Range("I107").Value = 813.46155
$ws.Range("J107").Value = 723.25
$ws.Range("K107").Value = 2440.38465
$ws.Range("L107").Value = 2169.75
$ws.Range("M107").Value = -520.38465
$ws.Range("N107").Value = -6009.75

# Row 132
$ws.Range("H132").Value = 6805303.5
$ws.Range("I132").Value = 1758.4642
$ws.Range("J132").Value = 15876697
$ws.Range("K132").Value = 5275.392599999999
$ws.Range("L132").Value = 47630091
$ws.Range("M132").Value = -2745.392599999999
$ws.Range("N132").Value = -47635151

# Row 136
$ws.Range("H136").Value = 3960.5356
$ws.Range("I136").Value = 2085.4707
$ws.Range("J136").Value = 6858.364
$ws.Range("K136").Value = 6256.4121
$ws.Range("L136").Value = 20575.092
$ws.Range("M136").Value = -3706.4121
$ws.Range("N136").Value = -25675.092
